$wb = $excel.ActiveWorkbook

# Work on the "Recommended but not categorized" sheet (4th tab).
$ws = $wb.Worksheets.Item("Recommended but not categorized")
$ws.Activate()

# Clear the title cell (A1) - it no longer holds the "Cory_Gaelic" label.
$ws.Range("A1").ClearContents()

# Delete row 2 (blank separator row) - everything below shifts up by one.
$ws.Rows(2).Delete()

# Delete the three blank separator rows that were at 68-70 (now at 67-69
# after the previous row-2 deletion shifted everything up by one).
$ws.Rows("67:69").Delete()

# Leave the selection on A4, matching the saved view state.
$ws.Range("A4").Select()
